$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 334, shifting rows 334:379 down to 335:380
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row 334 with its data
$ws.Cells.Item(334, 1).Value = 9
$ws.Cells.Item(334, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(334, 3).Value = "Metropolitana"
$ws.Cells.Item(334, 4).Value = 45127
$ws.Cells.Item(334, 5).Value = 13
$ws.Cells.Item(334, 6).Value = 100112001
$ws.Cells.Item(334, 7).Value = "Berenjena"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 70
$ws.Cells.Item(334, 11).Value = 7000
$ws.Cells.Item(334, 12).Value = 8000
$ws.Cells.Item(334, 13).Value = 7500
$ws.Cells.Item(334, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(334, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(334, 16).Value = 150
$ws.Cells.Item(334, 17).Value = 50
$ws.Cells.Item(334, 18).Value = "Hortaliza"

# Apply the date style (style index 2, numFmt 165) to D334, matching other date cells in column D
$ws.Cells.Item(334, 4).NumberFormat = $ws.Cells.Item(335, 4).NumberFormat
